$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.864.54"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.624.34"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'210.98"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'23.35"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").Value = "'0.256"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "'0.0879"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "1.618.00"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("D16").Value = "'65.28"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "27.857.13"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "'229.15"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0721"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'7.62"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D22").Value = "'4.31"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").Value = "'10.07"
$ws.Range("E23").Value = "  -5.91%  "
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("D26").Value = "'6.90"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "'0.0480"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").Value = "1.391.76"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  +11.62%  "
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "'0.0169"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "'0.857"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("D41").Value = "'1.03"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "'1.82"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").Value = "'65.58"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").Value = "1.765.79"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").Value = "'2.15"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("D48").Value = "'87.88"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("D50").Value = "'0.101"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("E51").Value = "  -0.48%  "
